$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Target cluster becomes "ECs", recompute numeric columns ---
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.528376666666667
$ws.Range("H2").Value = 4.58513
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.09433999999999999
$ws.Range("N2").Value = 0.28302
$ws.Range("O2").Value = 0.05191071108246543
$ws.Range("P2").Value = 0.05191071108246543
$ws.Range("Q2").Value = 0.1441870547333333
$ws.Range("R2").Value = 1.2976834926
$ws.Range("S2").Value = 0.05191071108246543
$ws.Range("T2").Value = 0.05191071108246543

# --- Row 3: Target cluster becomes "FAPs", recompute numeric columns ---
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.528376666666667
$ws.Range("H3").Value = 4.58513
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9431116666666667
$ws.Range("N3").Value = 2.829335
$ws.Range("O3").Value = 0.5189484550226392
$ws.Range("P3").Value = 0.5189484550226391
$ws.Range("Q3").Value = 1.441429865394445
$ws.Range("R3").Value = 12.97286878855
$ws.Range("S3").Value = 0.5189484550226392
$ws.Range("T3").Value = 0.5189484550226391

# --- Row 4: new row, Target cluster "sCs" ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf10"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.528376666666667
$ws.Range("H4").Value = 4.58513
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7798996666666667
$ws.Range("N4").Value = 2.339699
$ws.Range("O4").Value = 0.4291408338948954
$ws.Range("P4").Value = 0.4291408338948954
$ws.Range("Q4").Value = 1.191980452874444
$ws.Range("R4").Value = 10.72782407587
$ws.Range("S4").Value = 0.4291408338948954
$ws.Range("T4").Value = 0.4291408338948954
